$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.3191
$ws.Range("C4").Value = -12.9068
$ws.Range("A6").Value = -22.51490000000001
$ws.Range("A7").Value = -18.94449999999998
$ws.Range("A8").Value = -21.9018
$ws.Range("C8").Value = -12.82030000000001
$ws.Range("C9").Value = -10.1376
$ws.Range("C12").Value = -10.18479999999999
$ws.Range("A16").Value = -21.87920000000002
$ws.Range("C17").Value = -14.49549999999999
$ws.Range("C18").Value = -13.61939999999999
$ws.Range("C19").Value = -11.3029
$ws.Range("A20").Value = -19.88639999999999
$ws.Range("C20").Value = -12.0859
$ws.Range("A21").Value = -20.33369999999997
$ws.Range("C26").Value = -11.3784
$ws.Range("A28").Value = -21.9474
$ws.Range("A29").Value = -21.11509999999997
$ws.Range("A30").Value = -21.7265
$ws.Range("C31").Value = -12.43960000000001
$ws.Range("A32").Value = -21.2223
$ws.Range("C39").Value = -11.4345
$ws.Range("A40").Value = -20.42250000000001
$ws.Range("C40").Value = -12.34400000000001
$ws.Range("C41").Value = -11.98799999999999
$ws.Range("C42").Value = -11.6357
$ws.Range("C43").Value = -12.84129999999999
$ws.Range("A46").Value = -21.73380000000002
$ws.Range("C47").Value = -11.8893
$ws.Range("C48").Value = -11.59369999999999
$ws.Range("A51").Value = -21.67319999999999
$ws.Range("A52").Value = -22.29499999999999
$ws.Range("C54").Value = -13.0664
$ws.Range("A57").Value = -22.40040000000001
$ws.Range("A59").Value = -22.32970000000001
$ws.Range("A62").Value = -21.9309
$ws.Range("C62").Value = -12.86620000000001
$ws.Range("C63").Value = -10.3895
$ws.Range("C64").Value = -10.05529999999999
$ws.Range("A66").Value = -21.58520000000001
$ws.Range("A73").Value = -20.40899999999997
$ws.Range("A74").Value = -21.62429999999998
$ws.Range("C76").Value = -12.3793
$ws.Range("A77").Value = -20.21789999999999
$ws.Range("C81").Value = -13.92980000000001
$ws.Range("C84").Value = -13.4319
$ws.Range("C89").Value = -14.69139999999999
$ws.Range("A92").Value = -21.63020000000001
$ws.Range("C94").Value = -10.5637
$ws.Range("A100").Value = -21.96370000000002
